# Auto-generated COM-interop script.
# Adds the knaerot ('Knarot') section + references before the final
# sectPr, and updates the first-page header date.

$ErrorActionPreference = "Stop"

function Decode([string]$s) {
    $s = $s.Replace("@ae@", [string][char]228)
    $s = $s.Replace("@oe@", [string][char]246)
    $s = $s.Replace("@aa@", [string][char]229)
    $s = $s.Replace("@AE@", [string][char]196)
    $s = $s.Replace("@eacute@", [string][char]233)
    $s = $s.Replace("@ndash@", [string][char]8211)
    $s = $s.Replace("@lq@", [string][char]8220)
    $s = $s.Replace("@rq@", [string][char]8221)
    $s = $s.Replace("@sect@", [string][char]167)
    return $s
}

$d = $word.ActiveDocument

# --- locate the anchor paragraph ('BILAGA 1 - Fridlysta arter'), which is
# the last paragraph of the document body, immediately before the sectPr ---
$anchorIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs($anchorIndex)
if ($anchor.Range.Text -notmatch "BILAGA 1") {
    throw "Anchor paragraph not found where expected"
}

$insertionPoint = $anchor.Range
$insertionPoint.Collapse(0)  # wdCollapseEnd
$lastIndex = $anchorIndex

######################################################################
# Phase 1: create every paragraph and insert its PLAIN (unformatted)
# text. No Font.Italic assignment happens in this phase.
######################################################################

# paragraph 0
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Heading 1"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p0Start = $runRange.Start
$t0_0 = Decode "Kn@ae@rot @ndash@ ekologi samt krav p@aa@ livsmilj@oe@n"
$runRange.InsertAfter($t0_0)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 1
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p1Start = $runRange.Start
$t1_0 = Decode "Kn@ae@rot @ae@r fridlyst enligt 8 och 15 @sect@@sect@ artskyddsf@oe@rordningen och klassad som s@aa@rbar (VU) enligt r@oe@dlistan 2020. Kn@ae@rot @ae@r beroende av h@oe@g och j@ae@mn luftfuktighet i gamla, ost@oe@rda skogsmilj@oe@er och @ae@r k@ae@nslig f@oe@r snabba f@oe@r@ae@ndringar av ljus-/vindf@oe@rh@aa@llanden eller uttorkning. P@aa@ grund av ett alltf@oe@r intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 @aa@ren och i framtiden bed@oe@ms minskningstakten uppg@aa@ till 30 (20-40) %. Till f@oe@ljd av att arten har en dokumenterat h@oe@gre minskningstakt if@oe@rh@aa@llande till sin generationstid @ae@n vad som tidigare varit k@ae@nt (data fr@aa@n Riksskogstaxeringen) h@oe@jdes den till hotkategori s@aa@rbar (VU) i r@oe@dlistan 2020 (Artdatabanken, 2021)."
$runRange.InsertAfter($t1_0)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 2
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p2Start = $runRange.Start
$t2_0 = Decode "Samuel Johnsons doktorsavhandling "
$runRange.InsertAfter($t2_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t2_1 = Decode "@lq@Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation@lq@"
$r2_1Start = $runRange.Start
$runRange.InsertAfter($t2_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r2_1End = $runRange.Start
$t2_2 = Decode " (SLU, Uppsala 2014) visar att det kr@ae@vs v@ae@l tilltagna skyddszoner f@oe@r att kn@ae@rotens v@ae@xtplatser inte ska ta skada av skogsbruks@aa@tg@ae@rder i intilliggande omr@aa@den: "
$runRange.InsertAfter($t2_2)
$runRange.Collapse(0)  # wdCollapseEnd
$t2_3 = Decode "@lq@Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.@rq@ "
$r2_3Start = $runRange.Start
$runRange.InsertAfter($t2_3)
$runRange.Collapse(0)  # wdCollapseEnd
$r2_3End = $runRange.Start
$t2_4 = Decode "Vidare "
$runRange.InsertAfter($t2_4)
$runRange.Collapse(0)  # wdCollapseEnd
$t2_5 = Decode "@lq@More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).@rq@"
$r2_5Start = $runRange.Start
$runRange.InsertAfter($t2_5)
$runRange.Collapse(0)  # wdCollapseEnd
$r2_5End = $runRange.Start
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 3
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p3Start = $runRange.Start
$t3_0 = Decode "Johnsons (2014) rekommendation p@aa@ minst 50 meters breda skyddszoner runt kn@ae@rotens v@ae@xtplatser motsvarar en areal p@aa@ 0,78 hektar, vilket ligger i linje med andra studier som gjorts p@aa@ k@ae@nsliga skogsarter: "
$runRange.InsertAfter($t3_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t3_1 = Decode "@lq@In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).@rq@"
$r3_1Start = $runRange.Start
$runRange.InsertAfter($t3_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r3_1End = $runRange.Start
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 4
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p4Start = $runRange.Start
$t4_0 = Decode "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkid@eacute@n kn@ae@rots skyddsbehov. I uppsatsen ber@oe@rs problemet med uttorkning f@oe@r v@ae@xter, bl.a. f@oe@r kn@ae@rot, ett problem som blivit accentuerat p@aa@ grund av den p@aa@g@aa@ende klimatf@oe@r@ae@ndringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen unders@oe@ks omr@aa@den med tre olika avst@aa@nd fr@aa@n kalhyggeskant med avseende p@aa@ skydd bl.a. f@oe@r kn@ae@rot. Det f@oe@rsta omr@aa@det har avst@aa@nd upp till 20 m fr@aa@n hyggeskant (Strong edge effect), det andra 20 @ndash@ 40 m fr@aa@n hyggeskant (Weak edge effect) och det tredje avser st@oe@rre avst@aa@nd fr@aa@n hyggeskant, d@ae@r kanteffekten anses vara f@oe@rsumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt p@aa@ k@ae@nsliga och r@oe@dlistade skogsarter vid de kortare avst@aa@nden till hyggeskant, medan effekt av uttorkning inte konstaterades p@aa@ st@oe@rre avst@aa@nd (Interior). F@oe@r orkid@eacute@n kn@ae@rot fann man en rik f@oe@rekomst (upp till 0,06 dm2/m2) p@aa@ stort avst@aa@nd fr@aa@n hyggeskant (Interior), medan f@oe@rekomsten var liten eller n@ae@rmast f@oe@rsumbar i de omr@aa@den som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet p@aa@pekar att de allt oftare f@oe@rekommande torra somrarna ger ytterligare sk@ae@l att ut@oe@ka skyddsavst@aa@ndet fr@aa@n hyggen till den fuktkr@ae@vande arten kn@ae@rot (Koelmeijer m.fl., 2022)."
$runRange.InsertAfter($t4_0)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 5
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p5Start = $runRange.Start
$t5_0 = Decode "@AE@ven Skogsstyrelsens egen v@ae@gledning f@oe@r h@ae@nsyn till kn@ae@rot ligger i linje med ovanst@aa@ende forskningsstudier. Av v@ae@gledningen framg@aa@r det att f@oe@r med h@oe@g sannolikhet kunna bevara befintliga f@oe@rekomster kr@ae@vs relativt stora avs@ae@ttningar av uppvuxen skog med slutet och relativt t@ae@tt kronskikt. Som riktlinje kan kr@ae@vas ett avst@aa@nd p@aa@ 50 meter in fr@aa@n brynet f@oe@r att vidmakth@aa@lla ett fungerande mikroklimat. Detta inneb@ae@r att frist@aa@ende h@ae@nsynsytor f@oe@r m@aa@nga arter (k@ae@rlv@ae@xter, lavar och mossor) kan beh@oe@va ha en area @oe@verstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) f@oe@r att bibeh@aa@lla lokalklimatet. @AE@ven ganska sm@aa@ f@oe@r@ae@ndringar i form av f@oe@r@ae@ndrade ljus- och fuktighetsf@oe@rh@aa@llanden, till exempel till f@oe@ljd av gallring, kan leda till att arten f@oe@rsvinner till f@oe@ljd av konkurrens med mera ljuskr@ae@vande och snabbv@ae@xande arter (Skogsstyrelsen, 2022)."
$runRange.InsertAfter($t5_0)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 6
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Heading 2"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p6Start = $runRange.Start
$t6_0 = Decode "Referenser - kn@ae@rot"
$runRange.InsertAfter($t6_0)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 7
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p7Start = $runRange.Start
$t7_0 = Decode "de Graaf M & Roberts M.R., 2009. "
$runRange.InsertAfter($t7_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t7_1 = Decode "Short-term response of the herbaceous layer within leave patches after harvest. "
$r7_1Start = $runRange.Start
$runRange.InsertAfter($t7_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r7_1End = $runRange.Start
$t7_2 = Decode "Forest Ecology and Management 257, 1014-1025"
$runRange.InsertAfter($t7_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 8
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p8Start = $runRange.Start
$t8_0 = Decode "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. "
$runRange.InsertAfter($t8_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t8_1 = Decode "Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. "
$r8_1Start = $runRange.Start
$runRange.InsertAfter($t8_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r8_1End = $runRange.Start
$t8_2 = Decode "Ecological Applications, 22, 2049-2064 "
$runRange.InsertAfter($t8_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 9
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p9Start = $runRange.Start
$t9_0 = Decode "Koelmeijer, I. A., Ehrl@eacute@n, J., J@oe@nsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. "
$runRange.InsertAfter($t9_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t9_1 = Decode "Interactive effects of drought and edge exposure on old-growth forest understory species. "
$r9_1Start = $runRange.Start
$runRange.InsertAfter($t9_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r9_1End = $runRange.Start
$t9_2 = Decode "Landscape Ecology, 37, sid 1839-1853"
$runRange.InsertAfter($t9_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 10
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p10Start = $runRange.Start
$t10_0 = Decode "Rudolphi, J., J@oe@nsson, M. T., & Gustafsson, L., 2014. "
$runRange.InsertAfter($t10_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t10_1 = Decode "Biological legacies buffer local species extinction after logging. "
$r10_1Start = $runRange.Start
$runRange.InsertAfter($t10_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r10_1End = $runRange.Start
$t10_2 = Decode "Journal of Applied Ecology. 51, 53-62."
$runRange.InsertAfter($t10_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 11
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p11Start = $runRange.Start
$t11_0 = Decode "Skogsstyrelsen, 2022. "
$runRange.InsertAfter($t11_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t11_1 = Decode "V@ae@gledning f@oe@r h@ae@nsyn till kn@ae@rot. "
$r11_1Start = $runRange.Start
$runRange.InsertAfter($t11_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r11_1End = $runRange.Start
$t11_2 = Decode "https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/"
$runRange.InsertAfter($t11_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

# paragraph 12
$insertionPoint.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$newPara = $d.Paragraphs($lastIndex)
$newPara.Style = "Normal"
$runRange = $newPara.Range
$runRange.Collapse(1)  # wdCollapseStart
$p12Start = $runRange.Start
$t12_0 = Decode "SLU Artdatabanken, 2021. "
$runRange.InsertAfter($t12_0)
$runRange.Collapse(0)  # wdCollapseEnd
$t12_1 = Decode "Artfaktablad. Naturv@aa@rd @ndash@ artfakta. "
$r12_1Start = $runRange.Start
$runRange.InsertAfter($t12_1)
$runRange.Collapse(0)  # wdCollapseEnd
$r12_1End = $runRange.Start
$t12_2 = Decode "SLU Artdatabanken, Uppsala "
$runRange.InsertAfter($t12_2)
$runRange.Collapse(0)  # wdCollapseEnd
$insertionPoint = $d.Paragraphs($lastIndex).Range
$insertionPoint.Collapse(0)  # wdCollapseEnd

######################################################################
# Phase 2: now that every paragraph already exists, go back and apply
# italic formatting to the recorded sub-ranges. Nothing here sits at
# a live paragraph-insertion point, so no formatting bleeds across a
# paragraph boundary.
######################################################################

$d.Range($r2_1Start, $r2_1End).Font.Italic = $true
$d.Range($r2_3Start, $r2_3End).Font.Italic = $true
$d.Range($r2_5Start, $r2_5End).Font.Italic = $true
$d.Range($r3_1Start, $r3_1End).Font.Italic = $true
$d.Range($r7_1Start, $r7_1End).Font.Italic = $true
$d.Range($r8_1Start, $r8_1End).Font.Italic = $true
$d.Range($r9_1Start, $r9_1End).Font.Italic = $true
$d.Range($r10_1Start, $r10_1End).Font.Italic = $true
$d.Range($r11_1Start, $r11_1End).Font.Italic = $true
$d.Range($r12_1Start, $r12_1End).Font.Italic = $true

######################################################################
# Update the date in the first-page header
######################################################################
$sec = $d.Sections(1)
$fph = $sec.Headers(2)  # wdHeaderFooterFirstPage
$found = $fph.Range.Find.Execute(
    "2023-09-13", $false, $false, $false, $false, $false,
    $true, 1, $false, "2023-09-15", 2
)
if (-not $found) {
    throw "Date replacement in header failed"
}

Write-Host "Done."
